# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6926
$ws1.Range("F5").Value = 57
$ws1.Range("F6").Value = 1075
$ws1.Range("F7").Value = 162
$ws1.Range("F8").Value = 8

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6926
$ws4.Range("F5").Value = 57
$ws4.Range("F6").Value = 1075
$ws4.Range("F7").Value = 162
$ws4.Range("F9").Value = 8
